$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("D29").Value = ""
